$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25, shifting existing rows 25+ down by one.
$ws.Rows(25).Insert()

# Fill in the new row 25 with data (same constant columns as the rest of the sheet).
$ws.Range("A25").Value = 3
$ws.Range("B25").Value = "Femacal de La Calera"
$ws.Range("C25").Value = "Coquimbo"
$ws.Range("D25").Value = 44620
$ws.Range("E25").Value = 5
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100107
$ws.Range("H25").Value = "Otros"
$ws.Range("I25").Value = 100107011
$ws.Range("J25").Value = "Tuna"
$ws.Range("K25").Value = "Sin especificar"
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 70
$ws.Range("N25").Value = 13000
$ws.Range("O25").Value = 13000
$ws.Range("P25").Value = 13000
$ws.Range("Q25").Value = "$/caja 16 kilos"
$ws.Range("R25").Value = "Cabildo"
$ws.Range("S25").Value = 812
$ws.Range("T25").Value = 16
